# Update the CABLE ACERO Y ALAMBRE ROPA price list:
#  - A1 holds the price-list date (serial date value), bump it by one month
#    (45406 = 2024-04-24  ->  45436 = 2024-05-24)
#  - D22 holds the price per meter for the "ALAMBRE FORRADO" item (A-060)
#  - D37 holds the price per meter for the "CABLE DE ACERO" item (C-060)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A1").Value = 45436
$ws.Range("D22").Value = 155
$ws.Range("D37").Value = 206
